$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2.0
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.5439096666666666
$ws.Cells.Item(2, 8).Value = 1.631729
$ws.Cells.Item(2, 9).Value = 0.003493229883501837
$ws.Cells.Item(2, 10).Value = 0.003493229883501837
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 0.4702473333333333
$ws.Cells.Item(2, 14).Value = 1.410742
$ws.Cells.Item(2, 15).Value = 0.00903492226842282
$ws.Cells.Item(2, 16).Value = 0.00903492226842282
$ws.Cells.Item(2, 17).Value = 0.2557720703242222
$ws.Cells.Item(2, 18).Value = 2.301948632918
$ws.Cells.Item(2, 19).Value = 0.0000315610604631708
$ws.Cells.Item(2, 20).Value = 0.0000315610604631708
# Row 3
$ws.Cells.Item(3, 5).Value = 2.0
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.5439096666666666
$ws.Cells.Item(3, 8).Value = 1.631729
$ws.Cells.Item(3, 9).Value = 0.003493229883501837
$ws.Cells.Item(3, 10).Value = 0.003493229883501837
$ws.Cells.Item(3, 14).Value = 0.9584440000000001
$ws.Cells.Item(3, 15).Value = 0.006138235792679485
$ws.Cells.Item(3, 16).Value = 0.006138235792679485
$ws.Cells.Item(3, 17).Value = 0.1737689855195556
$ws.Cells.Item(3, 18).Value = 1.563920869676
$ws.Cells.Item(3, 19).Value = 0.00002144226870296857
$ws.Cells.Item(3, 20).Value = 0.00002144226870296857
# Row 4
$ws.Cells.Item(4, 5).Value = 2.0
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.5439096666666666
$ws.Cells.Item(4, 8).Value = 1.631729
$ws.Cells.Item(4, 9).Value = 0.003493229883501837
$ws.Cells.Item(4, 10).Value = 0.003493229883501837
$ws.Cells.Item(4, 13).Value = 1.047307
$ws.Cells.Item(4, 14).Value = 3.141921
$ws.Cells.Item(4, 15).Value = 0.02012204358311108
$ws.Cells.Item(4, 16).Value = 0.02012204358311108
$ws.Cells.Item(4, 17).Value = 0.5696404012676666
$ws.Cells.Item(4, 18).Value = 5.126763611409
$ws.Cells.Item(4, 19).Value = 0.00007029092396165002
$ws.Cells.Item(4, 20).Value = 0.00007029092396165002
# Row 5
$ws.Cells.Item(5, 5).Value = 2.0
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.5439096666666666
$ws.Cells.Item(5, 8).Value = 1.631729
$ws.Cells.Item(5, 9).Value = 0.003493229883501837
$ws.Cells.Item(5, 10).Value = 0.003493229883501837
$ws.Cells.Item(5, 13).Value = 50.21070966666667
$ws.Cells.Item(5, 14).Value = 150.632129
$ws.Cells.Item(5, 15).Value = 0.9647047983557866
$ws.Cells.Item(5, 16).Value = 0.9647047983557866
$ws.Cells.Item(5, 17).Value = 27.31009035789345
$ws.Cells.Item(5, 18).Value = 245.790813221041
$ws.Cells.Item(5, 19).Value = 0.003369935630374048
$ws.Cells.Item(5, 20).Value = 0.003369935630374048
# Row 6
$ws.Cells.Item(6, 9).Value = 0.00653284034046588
$ws.Cells.Item(6, 10).Value = 0.006532840340465881
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 0.4702473333333333
$ws.Cells.Item(6, 14).Value = 1.410742
$ws.Cells.Item(6, 15).Value = 0.00903492226842282
$ws.Cells.Item(6, 16).Value = 0.00903492226842282
$ws.Cells.Item(6, 17).Value = 0.4783304147459999
$ws.Cells.Item(6, 18).Value = 4.304973732714
$ws.Cells.Item(6, 19).Value = 0.0000590237046681261
$ws.Cells.Item(6, 20).Value = 0.0000590237046681261
# Row 7
$ws.Cells.Item(7, 9).Value = 0.00653284034046588
$ws.Cells.Item(7, 10).Value = 0.006532840340465881
$ws.Cells.Item(7, 14).Value = 0.9584440000000001
$ws.Cells.Item(7, 15).Value = 0.006138235792679485
$ws.Cells.Item(7, 16).Value = 0.006138235792679485
$ws.Cells.Item(7, 17).Value = 0.324972897972
$ws.Cells.Item(7, 18).Value = 2.924756081748
$ws.Cells.Item(7, 19).Value = 0.0000401001144057081
$ws.Cells.Item(7, 20).Value = 0.0000401001144057081
# Row 8
$ws.Cells.Item(8, 9).Value = 0.00653284034046588
$ws.Cells.Item(8, 10).Value = 0.006532840340465881
$ws.Cells.Item(8, 13).Value = 1.047307
$ws.Cells.Item(8, 14).Value = 3.141921
$ws.Cells.Item(8, 15).Value = 0.02012204358311108
$ws.Cells.Item(8, 16).Value = 0.02012204358311108
$ws.Cells.Item(8, 17).Value = 1.065309160023
$ws.Cells.Item(8, 18).Value = 9.587782440207
$ws.Cells.Item(8, 19).Value = 0.0001314540980523607
$ws.Cells.Item(8, 20).Value = 0.0001314540980523607
# Row 9
$ws.Cells.Item(9, 9).Value = 0.00653284034046588
$ws.Cells.Item(9, 10).Value = 0.006532840340465881
$ws.Cells.Item(9, 13).Value = 50.21070966666667
$ws.Cells.Item(9, 14).Value = 150.632129
$ws.Cells.Item(9, 15).Value = 0.9647047983557866
$ws.Cells.Item(9, 16).Value = 0.9647047983557866
$ws.Cells.Item(9, 17).Value = 51.073781555127
$ws.Cells.Item(9, 18).Value = 459.6640339961431
$ws.Cells.Item(9, 19).Value = 0.006302262423339685
$ws.Cells.Item(9, 20).Value = 0.006302262423339686
# Row 10
$ws.Cells.Item(10, 5).Value = 2.0
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.21601
$ws.Cells.Item(10, 8).Value = 0.6480299999999999
$ws.Cells.Item(10, 9).Value = 0.001387312330298533
$ws.Cells.Item(10, 10).Value = 0.001387312330298533
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 0.4702473333333333
$ws.Cells.Item(10, 14).Value = 1.410742
$ws.Cells.Item(10, 15).Value = 0.00903492226842282
$ws.Cells.Item(10, 16).Value = 0.00903492226842282
$ws.Cells.Item(10, 17).Value = 0.1015781264733333
$ws.Cells.Item(10, 18).Value = 0.9142031382599998
$ws.Cells.Item(10, 19).Value = 0.00001253425906627177
$ws.Cells.Item(10, 20).Value = 0.00001253425906627177
# Row 11
$ws.Cells.Item(11, 5).Value = 2.0
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.21601
$ws.Cells.Item(11, 8).Value = 0.6480299999999999
$ws.Cells.Item(11, 9).Value = 0.001387312330298533
$ws.Cells.Item(11, 10).Value = 0.001387312330298533
$ws.Cells.Item(11, 14).Value = 0.9584440000000001
$ws.Cells.Item(11, 15).Value = 0.006138235792679485
$ws.Cells.Item(11, 16).Value = 0.006138235792679485
$ws.Cells.Item(11, 17).Value = 0.06901116281333332
$ws.Cells.Item(11, 18).Value = 0.62110046532
$ws.Cells.Item(11, 19).Value = 0.00000851565020146404
$ws.Cells.Item(11, 20).Value = 0.00000851565020146404
# Row 12
$ws.Cells.Item(12, 5).Value = 2.0
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.21601
$ws.Cells.Item(12, 8).Value = 0.6480299999999999
$ws.Cells.Item(12, 9).Value = 0.001387312330298533
$ws.Cells.Item(12, 10).Value = 0.001387312330298533
$ws.Cells.Item(12, 13).Value = 1.047307
$ws.Cells.Item(12, 14).Value = 3.141921
$ws.Cells.Item(12, 15).Value = 0.02012204358311108
$ws.Cells.Item(12, 16).Value = 0.02012204358311108
$ws.Cells.Item(12, 17).Value = 0.2262287850699999
$ws.Cells.Item(12, 18).Value = 2.036059065629999
$ws.Cells.Item(12, 19).Value = 0.00002791555917365449
$ws.Cells.Item(12, 20).Value = 0.00002791555917365449
# Row 13
$ws.Cells.Item(13, 5).Value = 2.0
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.21601
$ws.Cells.Item(13, 8).Value = 0.6480299999999999
$ws.Cells.Item(13, 9).Value = 0.001387312330298533
$ws.Cells.Item(13, 10).Value = 0.001387312330298533
$ws.Cells.Item(13, 13).Value = 50.21070966666667
$ws.Cells.Item(13, 14).Value = 150.632129
$ws.Cells.Item(13, 15).Value = 0.9647047983557866
$ws.Cells.Item(13, 16).Value = 0.9647047983557866
$ws.Cells.Item(13, 17).Value = 10.84601539509667
$ws.Cells.Item(13, 18).Value = 97.61413855587
$ws.Cells.Item(13, 19).Value = 0.001338346861857143
$ws.Cells.Item(13, 20).Value = 0.001338346861857143
# Row 14
$ws.Cells.Item(14, 7).Value = 153.9268343333333
$ws.Cells.Item(14, 8).Value = 461.780503
$ws.Cells.Item(14, 9).Value = 0.9885866174457337
$ws.Cells.Item(14, 10).Value = 0.9885866174457337
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 12).Value = 1.0
$ws.Cells.Item(14, 13).Value = 0.4702473333333333
$ws.Cells.Item(14, 14).Value = 1.410742
$ws.Cells.Item(14, 15).Value = 0.00903492226842282
$ws.Cells.Item(14, 16).Value = 0.00903492226842282
$ws.Cells.Item(14, 17).Value = 72.38368337369177
$ws.Cells.Item(14, 18).Value = 651.4531503632259
$ws.Cells.Item(14, 19).Value = 0.008931803244225251
$ws.Cells.Item(14, 20).Value = 0.008931803244225251
# Row 15
$ws.Cells.Item(15, 7).Value = 153.9268343333333
$ws.Cells.Item(15, 8).Value = 461.780503
$ws.Cells.Item(15, 9).Value = 0.9885866174457337
$ws.Cells.Item(15, 10).Value = 0.9885866174457337
$ws.Cells.Item(15, 14).Value = 0.9584440000000001
$ws.Cells.Item(15, 15).Value = 0.006138235792679485
$ws.Cells.Item(15, 16).Value = 0.006138235792679485
$ws.Cells.Item(15, 17).Value = 49.17675026859244
$ws.Cells.Item(15, 18).Value = 442.590752417332
$ws.Cells.Item(15, 19).Value = 0.006068177759369344
$ws.Cells.Item(15, 20).Value = 0.006068177759369344
# Row 16
$ws.Cells.Item(16, 7).Value = 153.9268343333333
$ws.Cells.Item(16, 8).Value = 461.780503
$ws.Cells.Item(16, 9).Value = 0.9885866174457337
$ws.Cells.Item(16, 10).Value = 0.9885866174457337
$ws.Cells.Item(16, 13).Value = 1.047307
$ws.Cells.Item(16, 14).Value = 3.141921
$ws.Cells.Item(16, 15).Value = 0.02012204358311108
$ws.Cells.Item(16, 16).Value = 0.02012204358311108
$ws.Cells.Item(16, 17).Value = 161.2086510851403
$ws.Cells.Item(16, 18).Value = 1450.877859766263
$ws.Cells.Item(16, 19).Value = 0.01989238300192342
$ws.Cells.Item(16, 20).Value = 0.01989238300192342
# Row 17
$ws.Cells.Item(17, 7).Value = 153.9268343333333
$ws.Cells.Item(17, 8).Value = 461.780503
$ws.Cells.Item(17, 9).Value = 0.9885866174457337
$ws.Cells.Item(17, 10).Value = 0.9885866174457337
$ws.Cells.Item(17, 13).Value = 50.21070966666667
$ws.Cells.Item(17, 14).Value = 150.632129
$ws.Cells.Item(17, 15).Value = 0.9647047983557866
$ws.Cells.Item(17, 16).Value = 0.9647047983557866
$ws.Cells.Item(17, 17).Value = 7728.775588620098
$ws.Cells.Item(17, 18).Value = 69558.98029758089
$ws.Cells.Item(17, 19).Value = 0.9536942534402157
$ws.Cells.Item(17, 20).Value = 0.9536942534402157
